$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B4 switches from an inline text "5" to a real numeric value 5
$ws.Range("B4").Value = 5

# New row 5 with the latest form submission
$ws.Range("A5").Value = "TestGB"

# B5 keeps "25" as text (matches the source diff), so force text format
# before assigning -- otherwise Excel auto-converts the numeric-looking
# string to a number. Reset the style back to Normal afterwards so no
# stray number-format style lingers on the cell.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "25"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").Value = "hi"
